$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E2").Value = "5,50"
$ws.Range("F2").Value = "17 de jun., 18:34 UTC ·"

$ws.Range("E3").Value = "6,32"
$ws.Range("F3").Value = "17 de jun., 18:33 UTC ·"

$ws.Range("F4").Value = "17 de jun., 18:34 UTC ·"

$ws.Range("E5").Value = "7,39"
$ws.Range("F5").Value = "17 de jun., 18:32 UTC ·"

$ws.Range("F6").Value = "17 de jun., 18:34 UTC ·"

$ws.Range("E7").Value = "6,74"
$ws.Range("F7").Value = "17 de jun., 18:33 UTC ·"

$ws.Range("E8").Value = "3,56"
$ws.Range("F8").Value = "17 de jun., 18:33 UTC ·"

$ws.Range("F9").Value = "17 de jun., 18:33 UTC ·"

$ws.Range("E10").Value = "4,03"
$ws.Range("F10").Value = "17 de jun., 18:33 UTC ·"

$ws.Range("F11").Value = "17 de jun., 18:33 UTC ·"

$ws.Range("F12").Value = "17 de jun., 18:34 UTC ·"

$ws.Range("F13").Value = "17 de jun., 18:34 UTC ·"

$ws.Range("F14").Value = "17 de jun., 11:19 UTC ·"

$ws.Range("E15").Value = "0,0047"
$ws.Range("F15").Value = "17 de jun., 18:34 UTC ·"

$ws.Range("F16").Value = "17 de jun., 18:33 UTC ·"

$ws.Range("E17").Value = "1,10"
$ws.Range("F17").Value = "17 de jun., 18:33 UTC ·"

$ws.Range("E18").Value = "1,47"
$ws.Range("F18").Value = "17 de jun., 18:34 UTC ·"

$ws.Range("E19").Value = "4,28"
$ws.Range("F19").Value = "17 de jun., 18:33 UTC ·"

$ws.Range("F20").Value = "17 de jun., 18:34 UTC ·"

$ws.Range("E21").Value = "4,13"
$ws.Range("F21").Value = "17 de jun., 18:34 UTC ·"

$ws.Range("F22").Value = "17 de jun., 18:34 UTC ·"

$ws.Range("E23").Value = "3,48"
$ws.Range("F23").Value = "17 de jun., 18:33 UTC ·"

$ws.Range("F24").Value = "17 de jun., 18:34 UTC ·"

$ws.Range("F25").Value = "17 de jun., 18:34 UTC ·"
